$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data (A:E -> B:F)
$ws.Columns.Item(1).Insert()

# Copy the header formatting (bold, centered, bordered) from the
# neighbouring header cell onto the new header cell, then set its text.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false
$ws.Range("A1").Value2 = "ID"

# Fill in the new ID column values for each data row.
$ws.Cells.Item(2, 1).Value2 = "Hb 2"
$ws.Cells.Item(3, 1).Value2 = "Hb 3"
$ws.Cells.Item(4, 1).Value2 = "S 24"
$ws.Cells.Item(5, 1).Value2 = "S 28"
$ws.Cells.Item(6, 1).Value2 = "Hb 107"
$ws.Cells.Item(7, 1).Value2 = "Hb 66"
$ws.Cells.Item(8, 1).Value2 = "Hb 69"
$ws.Cells.Item(9, 1).Value2 = "Hb 95"
$ws.Cells.Item(10, 1).Value2 = "Hb 99"
$ws.Cells.Item(11, 1).Value2 = "Hb 92"
$ws.Cells.Item(12, 1).Value2 = "Hb 40"
$ws.Cells.Item(13, 1).Value2 = "Hb 41"
$ws.Cells.Item(14, 1).Value2 = "S 11"
$ws.Cells.Item(15, 1).Value2 = "Hb 57"
$ws.Cells.Item(16, 1).Value2 = "S 21"
$ws.Cells.Item(17, 1).Value2 = "S 22"
$ws.Cells.Item(18, 1).Value2 = "S 3"
$ws.Cells.Item(19, 1).Value2 = "S 4"
$ws.Cells.Item(20, 1).Value2 = "S 5"
$ws.Cells.Item(21, 1).Value2 = "Hb 74"
$ws.Cells.Item(22, 1).Value2 = "Hb 79"
$ws.Cells.Item(23, 1).Value2 = "Hb 32"
$ws.Cells.Item(24, 1).Value2 = "S 15"
$ws.Cells.Item(25, 1).Value2 = "S 16"
